$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 741
$ws1.Range("F6").Value = 2302
$ws1.Range("F8").Value = 1726
$ws1.Range("F9").Value = 2906
$ws1.Range("F10").Value = 165
$ws1.Range("F11").Value = 4342
$ws1.Range("F12").Value = 373
$ws1.Range("F13").Value = 201
$ws1.Range("F15").Value = 549
$ws1.Range("F16").Value = 260
$ws1.Range("F17").Value = 7
$ws1.Range("F18").Value = 126
$ws1.Range("F19").Value = 81
$ws1.Range("F20").Value = 100
$ws1.Range("F21").Value = 301
$ws1.Range("F22").Value = 4345
$ws1.Range("F24").Value = 3637
$ws1.Range("F25").Value = 1129
$ws1.Range("F26").Value = 210
$ws1.Range("F27").Value = 545
$ws1.Range("F28").Value = 4367
$ws1.Range("F30").Value = 512
$ws1.Range("F31").Value = 529
$ws1.Range("F32").Value = 483

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1032
$ws3.Range("F4").Value = 5

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1032
$ws4.Range("F5").Value = 5
$ws4.Range("F8").Value = 741
$ws4.Range("F9").Value = 2302
$ws4.Range("F11").Value = 1726
$ws4.Range("F13").Value = 2906
$ws4.Range("F14").Value = 165
$ws4.Range("F15").Value = 4342
$ws4.Range("F16").Value = 373
$ws4.Range("F17").Value = 201
$ws4.Range("F19").Value = 549
$ws4.Range("F20").Value = 260
$ws4.Range("F21").Value = 7
$ws4.Range("F22").Value = 126
$ws4.Range("F24").Value = 81
$ws4.Range("F25").Value = 100
$ws4.Range("F26").Value = 301
$ws4.Range("F27").Value = 4345
$ws4.Range("F29").Value = 3637
$ws4.Range("F30").Value = 1129
$ws4.Range("F31").Value = 210
$ws4.Range("F32").Value = 545
$ws4.Range("F33").Value = 4367
$ws4.Range("F35").Value = 512
$ws4.Range("F36").Value = 529
$ws4.Range("F37").Value = 483
